$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update translated text in column C (the 5 strings that were retranslated/reworded) ---
$ws.Range("C3").Value  = "Private Address"
$ws.Range("C5").Value  = "Operations"
$ws.Range("C6").Value  = "This virtual machine cannot be removed, and it must associate at least one security group."
$ws.Range("C8").Value  = "This container cannot be removed, and it must associate at least one security group."
$ws.Range("C10").Value = "Virtual Machine"

# --- Highlight the updated cells: red font + wrap text, no forced vertical centering ---
$updated = @("C3","C5","C6","C8","C10")
foreach ($addr in $updated) {
    $rng = $ws.Range($addr)
    $rng.Font.Color = 255
    $rng.WrapText = $true
    $rng.VerticalAlignment = -4107
}

# --- Row heights: rows 6 and 8 shrink from 60/45 to 30 ---
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 30

# --- Column C width tweak ---
$ws.Columns.Item(3).ColumnWidth = 37.92

# --- Selection moves to D8 ---
$ws.Range("D8").Select()

# --- Workbook calculation settings ---
$excel.Iteration = $false
$excel.CalculateBeforeSave = $false

# --- Normal style rename (best effort) ---
$wb.Styles.Item(1).NameLocal = "Normal"

# --- Page setup: portrait orientation ---
$ws.PageSetup().Orientation = 1
$ws.PageSetup().PaperSize = 0
